$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name / label text updates (shared-string reorder) ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Marzo de 2020 a las 19:42"
$ws.Range("A51").Value = "Catar"
$ws.Range("A52").Value = "Estonia"
$ws.Range("A77").Value = "Jordania"
$ws.Range("A78").Value = "San Marino"
$ws.Range("A79").Value = "Kuwait"
$ws.Range("A80").Value = "Costa Rica"
$ws.Range("A81").Value = "Republica de Macedonia"
$ws.Range("A82").Value = "Tunez"
$ws.Range("A83").Value = "Bosnia y Herzegovina"
$ws.Range("A84").Value = "Moldavia"
$ws.Range("A85").Value = "Albania"

# --- Numeric data updates ---
$ws.Range("B5").Value = 80589
$ws.Range("C5").Value = 6203
$ws.Range("E5").Value = 62013
$ws.Range("B6").Value = 79082
$ws.Range("C6").Value = 10871
$ws.Range("D6").Value = 1864
$ws.Range("E6").Value = 76075
$ws.Range("F6").Value = 2112
$ws.Range("G6").Value = 116
$ws.Range("H6").Value = 1143
$ws.Range("B10").Value = 29155
$ws.Range("C10").Value = 3922
$ws.Range("D10").Value = 4948
$ws.Range("E10").Value = 22511
$ws.Range("F10").Value = 3375
$ws.Range("G10").Value = 365
$ws.Range("H10").Value = 1696
$ws.Range("B12").Value = 11658
$ws.Range("C12").Value = 2129
$ws.Range("E12").Value = 10945
$ws.Range("G12").Value = 113
$ws.Range("H12").Value = 578
$ws.Range("F31").Value = 58
$ws.Range("B44").Value = 722
$ws.Range("C44").Value = 65
$ws.Range("E44").Value = 661
$ws.Range("B51").Value = 549
$ws.Range("C51").Value = 12
$ws.Range("D51").Value = 43
$ws.Range("E51").Value = 506
$ws.Range("H51").Value = 0
$ws.Range("B52").Value = 538
$ws.Range("C52").Value = 134
$ws.Range("D52").Value = 8
$ws.Range("E52").Value = 529
$ws.Range("H52").Value = 1
$ws.Range("B77").Value = 212
$ws.Range("C77").Value = 40
$ws.Range("D77").Value = 1
$ws.Range("E77").Value = 211
$ws.Range("F77").Value = 0
$ws.Range("H77").Value = 0
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 4
$ws.Range("E78").Value = 183
$ws.Range("F78").Value = 12
$ws.Range("H78").Value = 21
$ws.Range("B79").Value = 208
$ws.Range("C79").Value = 13
$ws.Range("D79").Value = 49
$ws.Range("E79").Value = 159
$ws.Range("F79").Value = 7
$ws.Range("H79").Value = 0
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 2
$ws.Range("E80").Value = 197
$ws.Range("F80").Value = 4
$ws.Range("H80").Value = 2
$ws.Range("B81").Value = 201
$ws.Range("C81").Value = 24
$ws.Range("D81").Value = 3
$ws.Range("E81").Value = 195
$ws.Range("F81").Value = 1
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 3
$ws.Range("B82").Value = 200
$ws.Range("C82").Value = 27
$ws.Range("E82").Value = 192
$ws.Range("F82").Value = 10
$ws.Range("G82").Value = 1
$ws.Range("H82").Value = 6
$ws.Range("B83").Value = 189
$ws.Range("C83").Value = 13
$ws.Range("E83").Value = 184
$ws.Range("F83").Value = 1
$ws.Range("H83").Value = 3
$ws.Range("B84").Value = 177
$ws.Range("D84").Value = 2
$ws.Range("E84").Value = 174
$ws.Range("F84").Value = 28
$ws.Range("G84").Value = 0
$ws.Range("H84").Value = 1
$ws.Range("B85").Value = 174
$ws.Range("C85").Value = 28
$ws.Range("D85").Value = 17
$ws.Range("E85").Value = 151
$ws.Range("F85").Value = 3
$ws.Range("G85").Value = 1
$ws.Range("H85").Value = 6
